$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new devlog entry row above the current top entry (row 9).
#    This shifts the existing rows 9..58 down to 10..59.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert(-4121)

# Re-apply the correct cell formatting to the new row 9 by copying the
# formats from row 10 (the row that used to be row 9 before the shift).
$ws.Range("A10:G10").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the row height for the new entry (it needs more vertical space).
$ws.Rows.Item(9).RowHeight = 90

# Fill in the content of the new devlog entry.
$ws.Range("B9").Value2 = "Milestone - Evaluation of prototypes has been completed and the function to evaluate prototypes is now called in main alongside all other unit tests. This means that all necessary code for the assignment has been written"
$ws.Range("C9").Value2 = 45777
$ws.Range("D9").Value2 = "read_data`ntime_many_to_many`nevaluate_database"
$ws.Range("E9").Value2 = "ChatGPT was used to write some of the read_data function, the lines written by ChatGPT are clearly marked in the function. The week 9 tutorial was used to learn about file reading and the week 4 tutorial code for timing code execution was used to time the many to many prototype function"
$ws.Range("F9").Value2 = "evaluate_database - Does not necessarily test the other functions implemented in this commit, but tests the prototypes and how they handle different m and n values."
$ws.Range("G9").Value2 = "The read_data function is taking a very long time to read all  2.5M records. This is to be expected, but a warning is still included in the evaluate_database function. The read_data function was modified many times to execute in the least amount of time possible."

# ---------------------------------------------------------------------------
# 2. Remove the now-duplicated blank spacer row (old row 20, pushed to row 21
#    by the insert above) that used to sit right under the devlog table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Delete()

# ---------------------------------------------------------------------------
# 3. Add a new blank row to the end of the trailing blank-row block.
# ---------------------------------------------------------------------------
$ws.Rows.Item(59).Insert(-4121)

# The row that used to be the final row (58) should now look like the other
# blank filler rows (21..57) instead of keeping its "last row" styling - that
# styling moves to the newly appended row 59.
$ws.Range("B57:G57").Copy()
$ws.Range("B58:G58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Update the active cell selection to match the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
